$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list refresh (coinranking snapshot moved from hour "16" to hour "17"
# on 2023-01-07) - updates Price (D), Volume 1h % (E) and Hora (G) columns for the
# affected rows. Values are written as literal text (matching the source sheet, where
# these columns are stored as text/inline strings) by forcing the cell NumberFormat to
# Text ("@") before assignment; this prevents Excel from auto-converting numeric-looking
# strings (e.g. "260.87", "0.86%", "17") into numbers/dates.
$updates = @(
    @{ Cell = "D2"; OldValue = "260.90"; NewValue = "260.87" },
    @{ Cell = "E2"; OldValue = "0.93%"; NewValue = "0.86%" },
    @{ Cell = "G2"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "D3"; OldValue = "27.05"; NewValue = "27.04" },
    @{ Cell = "E3"; OldValue = "0.65%"; NewValue = "0.53%" },
    @{ Cell = "G3"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "D4"; OldValue = "4.702"; NewValue = "4.703" },
    @{ Cell = "E4"; OldValue = "1.38%"; NewValue = "1.06%" },
    @{ Cell = "G4"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "D5"; OldValue = "0.06178"; NewValue = "0.06176" },
    @{ Cell = "E5"; OldValue = "3.69%"; NewValue = "3.58%" },
    @{ Cell = "G5"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "D6"; OldValue = "6.682"; NewValue = "6.679" },
    @{ Cell = "E6"; OldValue = "0.76%"; NewValue = "0.57%" },
    @{ Cell = "G6"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "D7"; OldValue = "0.8502"; NewValue = "0.8503" },
    @{ Cell = "G7"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "D8"; OldValue = "0.9138"; NewValue = "0.9105" },
    @{ Cell = "E8"; OldValue = "-1.06%"; NewValue = "-1.02%" },
    @{ Cell = "G8"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "D9"; OldValue = "0.1406"; NewValue = "0.1407" },
    @{ Cell = "E9"; OldValue = "1.34%"; NewValue = "1.65%" },
    @{ Cell = "G9"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "D10"; OldValue = "0.04651"; NewValue = "0.04649" },
    @{ Cell = "E10"; OldValue = "8.90%"; NewValue = "8.89%" },
    @{ Cell = "G10"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "D11"; OldValue = "0.07083"; NewValue = "0.07086" },
    @{ Cell = "E11"; OldValue = "0.91%"; NewValue = "1.01%" },
    @{ Cell = "G11"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "E12"; OldValue = "4.53%"; NewValue = "3.90%" },
    @{ Cell = "G12"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "D13"; OldValue = "0.09044"; NewValue = "0.09038" },
    @{ Cell = "E13"; OldValue = "-0.75%"; NewValue = "-0.82%" },
    @{ Cell = "G13"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "D14"; OldValue = "0.001539"; NewValue = "0.001533" },
    @{ Cell = "E14"; OldValue = "0.59%"; NewValue = "0.65%" },
    @{ Cell = "G14"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "D15"; OldValue = "0.0006167"; NewValue = "0.0006173" },
    @{ Cell = "E15"; OldValue = "1.94%"; NewValue = "1.95%" },
    @{ Cell = "G15"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "D16"; OldValue = "0.006048"; NewValue = "0.006047" },
    @{ Cell = "E16"; OldValue = "-0.63%"; NewValue = "-0.49%" },
    @{ Cell = "G16"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "D17"; OldValue = "3.453"; NewValue = "3.452" },
    @{ Cell = "E17"; OldValue = "0.03%"; NewValue = "0.01%" },
    @{ Cell = "G17"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "D18"; OldValue = "3.165"; NewValue = "3.167" },
    @{ Cell = "E18"; OldValue = "1.33%"; NewValue = "1.25%" },
    @{ Cell = "G18"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "E19"; OldValue = "1.18%"; NewValue = "1.17%" },
    @{ Cell = "G19"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "E20"; OldValue = "-0.88%"; NewValue = "-0.87%" },
    @{ Cell = "G20"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "E21"; OldValue = "0.22%"; NewValue = "0.21%" },
    @{ Cell = "G21"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "E22"; OldValue = "4.05%"; NewValue = "1.65%" },
    @{ Cell = "G22"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "D23"; OldValue = "0.04240"; NewValue = "0.04216" },
    @{ Cell = "E23"; OldValue = "0.38%"; NewValue = "-0.13%" },
    @{ Cell = "G23"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "D24"; OldValue = "0.001218"; NewValue = "0.001217" },
    @{ Cell = "E24"; OldValue = "0.08%"; NewValue = "0.09%" },
    @{ Cell = "G24"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "E25"; OldValue = "-5.65%"; NewValue = "-5.74%" },
    @{ Cell = "G25"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "E26"; OldValue = "0.11%"; NewValue = "0.14%" },
    @{ Cell = "G26"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "E27"; OldValue = "-7.80%"; NewValue = "-7.81%" },
    @{ Cell = "G27"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "G28"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "G29"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "G30"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "G31"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "G32"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "G33"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "G34"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "G35"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "G36"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "G37"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "G38"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "G39"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "E40"; OldValue = "1.49%"; NewValue = "1.44%" },
    @{ Cell = "G40"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "D41"; OldValue = "0.1112"; NewValue = "0.1111" },
    @{ Cell = "E41"; OldValue = "0.13%"; NewValue = "-0.14%" },
    @{ Cell = "G41"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "D42"; OldValue = "0.004090"; NewValue = "0.004098" },
    @{ Cell = "E42"; OldValue = "7.94%"; NewValue = "8.25%" },
    @{ Cell = "G42"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "E43"; OldValue = "9.28%"; NewValue = "9.26%" },
    @{ Cell = "G43"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "E44"; OldValue = "-10.05%"; NewValue = "-10.03%" },
    @{ Cell = "G44"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "D45"; OldValue = "0.00005162"; NewValue = "0.00005164" },
    @{ Cell = "E45"; OldValue = "0.15%"; NewValue = "0.32%" },
    @{ Cell = "G45"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "E46"; OldValue = "0.12%"; NewValue = "0.14%" },
    @{ Cell = "G46"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "G47"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "D48"; OldValue = "0.1682"; NewValue = "0.1667" },
    @{ Cell = "E48"; OldValue = "-23.92%"; NewValue = "-24.59%" },
    @{ Cell = "G48"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "E49"; OldValue = "0.12%"; NewValue = "0.14%" },
    @{ Cell = "G49"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "E50"; OldValue = "0.12%"; NewValue = "0.14%" },
    @{ Cell = "G50"; OldValue = "16"; NewValue = "17" },
    @{ Cell = "G51"; OldValue = "16"; NewValue = "17" }
)

foreach ($update in $updates) {
    $cell = $ws.Range($update.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $update.NewValue
}
